# Add MOU from Washington DC (DC BWC RCT study), and disambiguate the
# existing ASU "BWC RCT" study by renaming it to "ASU BWC RCT".
#
# Shared-string insertion order matters (Excel appends newly-seen unique
# strings to xl/sharedStrings.xml in the order they are first written), so
# the operations below are ordered to reproduce:
#   161 DC BWC RCT
#   162 A randomized control trial evaluating the effects of police body-worn cameras
#   163 Washington DC
#   164 https://doi.org/10.1073/pnas.1814773116
#   165 ASU BWC RCT
#   166 2023-FOIA-08968.pdf

$wb = $excel.ActiveWorkbook

$wsContracts = $wb.Worksheets.Item("Contracts")
$wsPapers    = $wb.Worksheets.Item("Papers")
$wsStudies   = $wb.Worksheets.Item("Studies")
$wsCities    = $wb.Worksheets.Item("Cities")

# --- 1) Studies: new row 21 -> "DC BWC RCT" / "Arnold Foundation" -----------
# First-ever use of "DC BWC RCT" -> becomes shared string 161.
$wsStudies.Activate()
$wsStudies.Range("A21").Value = "DC BWC RCT"
$wsStudies.Range("B21").Value = "Arnold Foundation"
$wsStudies.Range("B20:B20").Copy()
$wsStudies.Range("B21").PasteSpecial(-4122)
$wsStudies.Range("A7").Select()

# --- 2) Papers: new row 33 ---------------------------------------------------
# First-ever use of the paper title -> becomes shared string 162.
$wsPapers.Activate()
$wsPapers.Range("A33").Value = "A randomized control trial evaluating the effects of police body-worn cameras"
$wsPapers.Range("B33").Value = "DC BWC RCT"

# --- 3) Contracts: new row 44 (Contract With only, for now) -----------------
# First-ever use of "Washington DC" -> becomes shared string 163.
$wsContracts.Activate()
$wsContracts.Range("A44").Value = "Washington DC"
$wsContracts.Range("B44").Value = "DC BWC RCT"
$wsContracts.Range("A43:B43").Copy()
$wsContracts.Range("A44:B44").PasteSpecial(-4122)

# --- 4) Cities: new row 47 ---------------------------------------------------
$wsCities.Activate()
$wsCities.Range("A47").Value = "A randomized control trial evaluating the effects of police body-worn cameras"
$wsCities.Range("B47").Value = "Washington DC"
$wsCities.Range("B46:B46").Copy()
$wsCities.Range("B47").PasteSpecial(-4122)

# --- 5) Papers: finish row 33, add URL + hyperlink ---------------------------
# First-ever use of the DOI URL -> becomes shared string 164.
$wsPapers.Activate()
$wsPapers.Range("C33").Value = "https://doi.org/10.1073/pnas.1814773116"
$wsPapers.Hyperlinks.Add($wsPapers.Range("C33"), "https://doi.org/10.1073/pnas.1814773116")
$wsPapers.Range("A33").Select()

# --- 6) Contracts: finish row 44 with filename -------------------------------
# First-ever use of this filename -> becomes shared string 166 (after the
# "ASU BWC RCT" rename below claims 165).

# --- 7) Rename the existing ASU "BWC RCT" study to "ASU BWC RCT" ------------
# First rename -> introduces shared string 165.
$wsContracts.Activate()
$wsContracts.Range("B23").Value = "ASU BWC RCT"
$wsContracts.Range("B24").Value = "ASU BWC RCT"
$wsContracts.Range("B25").Value = "ASU BWC RCT"
$wsContracts.Range("B26").Value = "ASU BWC RCT"
$wsContracts.Range("B27").Value = "ASU BWC RCT"
$wsContracts.Range("B28").Value = "ASU BWC RCT"
$wsContracts.Range("B29").Value = "ASU BWC RCT"
$wsContracts.Range("B30").Value = "ASU BWC RCT"
$wsContracts.Range("B31").Value = "ASU BWC RCT"
$wsContracts.Range("B32").Value = "ASU BWC RCT"
$wsContracts.Range("B33").Value = "ASU BWC RCT"
$wsContracts.Range("B34").Value = "ASU BWC RCT"
$wsContracts.Range("B35").Value = "ASU BWC RCT"
$wsContracts.Range("B36").Value = "ASU BWC RCT"
$wsContracts.Range("B37").Value = "ASU BWC RCT"
$wsContracts.Range("B38").Value = "ASU BWC RCT"
$wsContracts.Range("B39").Value = "ASU BWC RCT"
$wsContracts.Range("B40").Value = "ASU BWC RCT"
$wsContracts.Range("B41").Value = "ASU BWC RCT"
# Note: row 43 ("Tempe" / "BWC RCT" / "Tempe ASU Services Agreement.pdf")
# intentionally keeps the original "BWC RCT" study name.

$wsPapers.Activate()
$wsPapers.Range("B6").Value = "ASU BWC RCT"
$wsPapers.Range("B7").Value = "ASU BWC RCT"
$wsPapers.Range("B8").Value = "ASU BWC RCT"
$wsPapers.Range("B9").Value = "ASU BWC RCT"
$wsPapers.Range("B10").Value = "ASU BWC RCT"
$wsPapers.Range("B11").Value = "ASU BWC RCT"
$wsPapers.Range("B12").Value = "ASU BWC RCT"

$wsStudies.Activate()
$wsStudies.Range("A6").Value = "ASU BWC RCT"

# --- 8) Contracts: finish row 44 with the filename --------------------------
# First-ever use of this filename -> becomes shared string 166.
$wsContracts.Activate()
$wsContracts.Range("C44").Value = "2023-FOIA-08968.pdf"

# --- Restore selections / active sheet to match the final view state -------
$wsCities.Activate()
$wsCities.Range("A48").Select()

$wsStudies.Activate()
$wsStudies.Range("A7").Select()

$wsPapers.Activate()
$wsPapers.Range("A33").Select()

$wsContracts.Activate()
$wsContracts.Range("A2:A3").Select()
